$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.515.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '''1.877.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''0.7190'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.57%  '
$ws.Range("D6").Value = '''242.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.98%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.07957'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").Value = '''0.3114'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.12%  '
$ws.Range("D10").Value = '''25.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.35%  '
$ws.Range("D11").Value = '''0.08275'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").Value = '''0.7324'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.78%  '
$ws.Range("D13").Value = '''5.297'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("D14").Value = '''1.867.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").Value = '''91.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '''29.502.87'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '''5.935'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.71%  '
$ws.Range("D18").Value = '''246.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.26%  '
$ws.Range("D19").Value = '''0.000007895'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("D20").Value = '''13.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").Value = '''2.117.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''8.029'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.82%  '
$ws.Range("D24").Value = '''1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '''0.1622'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +13.69%  '
$ws.Range("D26").Value = '''163.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.23%  '
$ws.Range("D27").Value = '''9.064'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("D28").Value = '''18.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = '''1.360'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.81%  '
$ws.Range("D30").Value = '''1.496'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("D31").Value = '''4.397'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("D32").Value = '''4.115'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.53%  '
$ws.Range("D33").Value = '''0.05275'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.36%  '
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("D35").Value = '''1.201'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.74%  '
$ws.Range("D36").Value = '''0.7277'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.69%  '
$ws.Range("D37").Value = '''2.683'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").Value = '''0.01875'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("D39").Value = '''1.208.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.23%  '
$ws.Range("D40").Value = '''2.704'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").Value = '''0.9103'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("D42").Value = '''73.88'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.32%  '
$ws.Range("D43").Value = '''6.142'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("D44").Value = '''1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = '''102.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.60%  '
$ws.Range("D46").Value = '''0.5292'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '''2.012.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").Value = '''1.800'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.56%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000121'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.76%  '
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").Value = '''2.942'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.01%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''9.382'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.32%  '
